# Added periodic & upfront related scenarios
# Update the "repaymentstrategy" value on the ProductLoanInput sheet from
# "Mifos style" to "Penalties, Fees, Interest, Principal order", applying a
# left/top aligned style, and move the active selection to that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$cell = $ws.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160
$cell.Select()
